$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Write the new bug text. The order here controls the shared-string table
#    insertion order, so it is deliberately NOT simple top-to-bottom row
#    order (it mirrors how the sheet was actually edited: the "Pagination
#    broken" / "Can't search" rows were filled in before the existing row 4
#    was rewritten into the "Deleted computer confirmation..." text).
# ---------------------------------------------------------------------------
$ws.Range("A5").Value  = "Can't search with `$ in the name"
$ws.Range("A6").Value  = "Pagination broken"
$ws.Range("B6").Value  = "1. search 'z'`n2. go to next page"
$ws.Range("C6").Value  = "displays 11 to 20 out of 14. Before clicking next it displays 1 - 10 of 24 "
$ws.Range("A4").Value  = "Deleted computer confirmation message doesn't display name"
$ws.Range("A7").Value  = "Able to set dates in future"
$ws.Range("A8").Value  = "Able to set dates  from far in the past"
$ws.Range("A9").Value  = "Date error -  ('yyyy-MM-dd') - Inconsistent case"
$ws.Range("A10").Value = "Cannot easily navigate to first and last pages of computer list"
$ws.Range("A11").Value = "Deleting a computer from search returns you to page 1 and clears search"
$ws.Range("A12").Value = "When reloading page or adding new computer occasionally the page styling isn't visible for a split second"
$ws.Range("A13").Value = "Errors when setting fields on add/edit computer don't inform you of the error"
$ws.Range("A14").Value = "Able to add computers where the discontinued date is before the introduced date."

# ---------------------------------------------------------------------------
# 2. Header cells (row 1 / row 3) - wrap text, values unchanged
# ---------------------------------------------------------------------------
$ws.Range("A1").WrapText = $true
$ws.Range("A3").WrapText = $true

# Seed the cellXfs cache in the same order the original author created it in
# (wrap-only, then vertical-top-only, then wrap+vertical-top) so the style
# indices line up: 1=wrap, 2=vertical-top, 3=wrap+vertical-top.
$ws.Range("B4").VerticalAlignment = -4160   # xlTop -> creates style index 2
$ws.Range("A4").WrapText = $true
$ws.Range("A4").VerticalAlignment = -4160   # -> creates style index 3

# ---------------------------------------------------------------------------
# 3. Apply per-cell alignment across the whole new block (A4:C14)
# ---------------------------------------------------------------------------
$rowHeights = @{ 4 = 30; 6 = 30; 10 = 30; 11 = 30; 12 = 45; 13 = 30; 14 = 30 }

for ($r = 4; $r -le 14; $r++) {
    $aCell = $ws.Cells.Item($r, 1)
    $bCell = $ws.Cells.Item($r, 2)
    $cCell = $ws.Cells.Item($r, 3)

    $aCell.WrapText = $true
    $aCell.VerticalAlignment = -4160

    $bCell.VerticalAlignment = -4160
    $cCell.VerticalAlignment = -4160

    if ($r -eq 6) {
        $bCell.WrapText = $true
    }

    if ($rowHeights.ContainsKey($r)) {
        $ws.Rows.Item($r).RowHeight = $rowHeights[$r]
    }
}

# ---------------------------------------------------------------------------
# 4. Column widths
# ---------------------------------------------------------------------------
$ws.Columns.Item(1).ColumnWidth = 49.66666666666666
$ws.Columns.Item(2).ColumnWidth = 32.5
$ws.Columns.Item(3).ColumnWidth = 62.33333333333333

# ---------------------------------------------------------------------------
# 5. Table resize + style
# ---------------------------------------------------------------------------
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A3:C14"))
$lo.TableStyle = "TableStyleMedium9"

# ---------------------------------------------------------------------------
# 6. Page setup
# ---------------------------------------------------------------------------
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

# ---------------------------------------------------------------------------
# 7. Selection
# ---------------------------------------------------------------------------
$ws.Range("B10").Select()
